# Replace placeholder "null" text in the states sheet's CITY column with
# the real numeric data, then restore the "states" tab as the active /
# selected sheet (with C3 selected) the way the workbook was left after
# this edit.

$wb = $excel.ActiveWorkbook

$states = $wb.Worksheets.Item("states")

# CITY column held literal text "null" for every data row - fill in the
# real numbers instead.
$states.Range("C2").Value = 503
$states.Range("C3").Value = 1873
$states.Range("C4").Value = 8068
$states.Range("C5").Value = 2185
$states.Range("C6").Value = 469
$states.Range("C7").Value = 2918

# The workbook was saved with the "states" sheet active (instead of
# "friends") and C3 selected on it.
$states.Activate()
$states.Range("C3").Select()
